$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header / unit / type cells for columns G and H ---
$ws.Range("G1").Value = "connectingTime"
$ws.Range("H1").Value = "passwordDigits"

$ws.Range("G2").Value = "해킹 서버 연결 시간(s)"
$ws.Range("H2").Value = "암호 자릿수"

$ws.Range("G3").Value = "int"
$ws.Range("H3").Value = "int"

$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 3

# --- Existing cells on row 2 that now say "미사용" instead of "sec" / blank ---
$ws.Range("B2").Value = "미사용"
$ws.Range("C2").Value = "미사용"
$ws.Range("D2").Value = "미사용"
$ws.Range("E2").Value = "미사용"

# --- Updated spy mission description ---
$ws.Range("F4").Value = "get 3 letters of the password before each three missions.`nInteract and wait the time for connecting to server.`nThen insert three letters you remembered."

# --- Copy cell formatting (border/fill/alignment) from existing styled cells ---
# Header-row style (gray fill, centered, bordered) -> row1 & row3 new cells
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null
$ws.Range("G3:H3").PasteSpecial(-4122) | Out-Null

# Plain bordered/centered style -> row2 & row4 new cells
$ws.Range("A2").Copy() | Out-Null
$ws.Range("G2:H2").PasteSpecial(-4122) | Out-Null
$ws.Range("G4:H4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Column widths for the new columns ---
# (ColumnWidth is quantized by the host to 1/7-character pixel steps, so the
#  assigned values are pre-compensated to land as close as possible to the
#  target stored widths of 25 and 18.375 characters.)
$ws.Columns.Item(7).ColumnWidth = 24.2857142857143
$ws.Columns.Item(8).ColumnWidth = 17.6607142857143

# --- Row 4 height shrinks now that the description is shorter ---
$ws.Rows.Item(4).RowHeight = 49.5

# --- Selection moves to F16 ---
$ws.Range("F16").Select() | Out-Null
